# Added 4wk low sales check:
# Updated MyForecast (D), Inventory Coverage (H), Reorder Urgency (J) and
# Seasonality Index (L) values on the "Forecast Comparison" sheet, and
# refreshed the dependent totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---

# Row 2 (W10)
$wsForecast.Range("D2").Value = 711
$wsForecast.Range("H2").Value = 8.16
$wsForecast.Range("L2").Value = 1

# Row 3 (W11)
$wsForecast.Range("D3").Value = 689
$wsForecast.Range("H3").Value = 7.39
$wsForecast.Range("L3").Value = 1.13

# Row 4 (W12)
$wsForecast.Range("D4").Value = 666
$wsForecast.Range("H4").Value = 6.61
$wsForecast.Range("L4").Value = 1.19

# Row 5 (W13)
$wsForecast.Range("D5").Value = 656
$wsForecast.Range("H5").Value = 5.69
$wsForecast.Range("L5").Value = 1.08

# Row 6 (W14)
$wsForecast.Range("D6").Value = 656
$wsForecast.Range("H6").Value = 4.69
$wsForecast.Range("L6").Value = 1.04

# Row 7 (W15)
$wsForecast.Range("D7").Value = 653
$wsForecast.Range("H7").Value = 3.71
$wsForecast.Range("L7").Value = 0.9

# Row 8 (W16)
$wsForecast.Range("D8").Value = 641
$wsForecast.Range("H8").Value = 2.76

# Row 9 (W17)
$wsForecast.Range("D9").Value = 630
$wsForecast.Range("H9").Value = 1.79
$wsForecast.Range("L9").Value = 0.9

# Row 10 (W18)
$wsForecast.Range("D10").Value = 635
$wsForecast.Range("H10").Value = 0.79
$wsForecast.Range("J10").Value = "Urgent"
$wsForecast.Range("L10").Value = 0.9399999999999999

# Row 11 (W19)
$wsForecast.Range("D11").Value = 654
$wsForecast.Range("H11").Value = 0

# Row 12 (W20)
$wsForecast.Range("D12").Value = 661
$wsForecast.Range("L12").Value = 0.9

# Row 13 (W21)
$wsForecast.Range("D13").Value = 643
$wsForecast.Range("L13").Value = 1.01

# Row 14 (W22)
$wsForecast.Range("D14").Value = 615
$wsForecast.Range("L14").Value = 0.93

# Row 15 (W23)
$wsForecast.Range("D15").Value = 608
$wsForecast.Range("L15").Value = 1.09

# Row 16 (W24)
$wsForecast.Range("D16").Value = 635
$wsForecast.Range("L16").Value = 1.05

# Row 17 (W25)
$wsForecast.Range("D17").Value = 669
$wsForecast.Range("L17").Value = 0.97

# --- Summary sheet ---

$wsSummary.Range("B9").Value = 10422
$wsSummary.Range("B10").Value = 5302
$wsSummary.Range("B11").Value = 2722
$wsSummary.Range("B12").Value = 711
$wsSummary.Range("B14").Value = 608
